# "aggiornamento diario di bordo"
# Update the "Riunioni interne di allineamento" bullet on slide 4 to reflect
# that alignment meetings are now planned via a shared calendar, and grow the
# text box that holds it to fit the longer line.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(14)   # "CasellaDiTesto 29"

$oldLine = "Riunioni interne di allineamento ogni due giorni"
$newLine = "Riunioni interne di allineamento (pianificate in un calendario condiviso)"

$tr = $sh.TextFrame.TextRange
$fullText = $tr.Text
$startPos = $fullText.IndexOf($oldLine)

if ($startPos -ge 0) {
    $target = $tr.Characters($startPos + 1, $oldLine.Length)
    $target.Text = $newLine
}

# Grow the text box height (EMU 4222438 -> 4822602) now that the bullet wraps
# onto more lines; width/position stay the same.
$sh.Height = 4822602 / 914400 * 72
